$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Drop the stray "_GoBack" bookmark that currently sits after
#    "Find the user and take the policies to call insurance policy function
#    for filtering." -- it will be re-created below at the new last-edit
#    location.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Locate the "FILTER FUNCTIONS (EXPIRY DATE, CAR)" note and strike
#    through just the word "CAR" (the rest of the run keeps its original
#    formatting), mirroring the other struck-through red notes already in
#    this document.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("FILTER FUNCTIONS (EXPIRY DATE, CAR)", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

$phraseStart = $rng.Start
$phraseEnd = $rng.End

$carRng = $d.Range($phraseStart, $phraseEnd)
[void]$carRng.Find.Execute("CAR", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$carRng.Font.StrikeThrough = $true

# ---------------------------------------------------------------------------
# 3) Re-insert the "_GoBack" bookmark right after the closing ")" (the new
#    last-edit position), collapsed to a single point. Adding a
#    zero-length bookmark directly at that offset trips an engine quirk, so
#    insert a throwaway character, bookmark the one-character range, then
#    delete the character again -- the bookmark correctly collapses back to
#    a point and survives the deletion.
# ---------------------------------------------------------------------------
$insPoint = $d.Range($phraseEnd, $phraseEnd)
$insPoint.InsertAfter("X")

$tmpRng = $d.Range($phraseEnd, $phraseEnd + 1)
$d.Bookmarks.Add("_GoBack", $tmpRng)

$tmpRng2 = $d.Range($phraseEnd, $phraseEnd + 1)
$tmpRng2.Text = ""
